$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. "Every since" -> "Ever since" (typo fix)
Replace-Text "Every since the advent" "Ever since the advent"

# 2. Split the combined "ethical and social responsibilities" quote into two separate quoted terms
Replace-Text "assume “ethical and social responsibilities” to their communities." "assume “ethical” and “social” responsibilities to their communities."

# 3. "scope that ethical" -> "scope than ethical"
Replace-Text "encompass a much broader scope that ethical ones" "encompass a much broader scope than ethical ones"

# 4. Remove "not only " before "asserts that corporations"
Replace-Text "This concept not only asserts that corporations" "This concept asserts that corporations"

# 5 & 6. "today's climate" -> "today's business climate"; "for to advertise" -> "for advertising"
Replace-Text "in today’s climate, it would make sense that organizations also use their home page for to advertise their philanthropy." "in today’s business climate, it would make sense that organizations also use their home page for advertising their philanthropy."

# 7. "for any organization." -> "for most organizations."
Replace-Text "I would make the following recommendations for any organization." "I would make the following recommendations for most organizations."

# 8. "Another Example" -> "Another example"
Replace-Text ". Another Example is that Google" ". Another example is that Google"

# 9. "values overlap those with their respective communities." -> "values that overlap values of their respective communities."
Replace-Text "business is based on values overlap those with their respective communities." "business is based on values that overlap values of their respective communities."

# 10. "Thermo Fisher, the larger organization" -> "Thermo Fisher Scientific, the larger organization"
Replace-Text "However, Thermo Fisher, the larger organization" "However, Thermo Fisher Scientific, the larger organization"

# 11. "occurrences of the how" -> "occurrences of how"
Replace-Text "Publishing occurrences of the how the organization donated" "Publishing occurrences of how the organization donated"

# 12. "are good public relations." -> "are good for public relations."
Replace-Text "Even though these actions may not save a struggling company, they are good public relations." "Even though these actions may not save a struggling company, they are good for public relations."

# 13. Reference list: "What we believe" -> "What We Believe" (italic citation, not the in-text quoted phrase).
# Restrict the search to the references paragraph so the body-text occurrence
# (inside quotation marks) is left untouched, and start the match one
# character inside the run (after the leading "W") so the italic run-level
# formatting boundary right before "What" is preserved by the replace engine.
$refPara = $d.Paragraphs(13).Range
$result = $refPara.Find.Execute("hat we believe", $true, $false, $false, $false, $false, $true, 1, $false, "hat We Believe", 2)
if (-not $result) {
    Write-Output "NOT FOUND: hat we believe (reference)"
}

# The "_GoBack" bookmark marks the last edit position. Since the final edit
# of this pass is the "What We Believe" citation fix above, move the
# bookmark there (right after "What We B", before "elieve"), matching where
# Word itself would leave it after typing that change.
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()
$refPara2 = $d.Paragraphs(13).Range
$found = $refPara2.Find.Execute("What We B", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($refPara2.End, $refPara2.End)
    $d.Bookmarks.Add("_GoBack", $target)
} else {
    Write-Output "NOT FOUND: What We B (bookmark anchor)"
}
